$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update quantities: OutTireLeft.stl and OutTireRight.stl counts go from 1 to 2
$ws.Range("C12").Value = 2
$ws.Range("C13").Value = 2

# Scroll the sheet view back to the top-left (removes the stashed
# topLeftCell="A22") and move the selection to I18
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("I18").Select()
